$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values like "0.140" or "592.71" that Excel
# would otherwise auto-parse into numbers (dropping trailing zeros / using
# scientific notation). Force text storage, then drop the temporary "@"
# number format again so the cell keeps its original (unformatted) style.
function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "62.483.54"
$ws.Range("E2").Value = "  +0.09%  "
Set-TextValue "D3" "3.125.34"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "592.71"
$ws.Range("E5").Value = "  +1.41%  "
Set-TextValue "D6" "133.36"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue "D8" "3.127.84"
$ws.Range("E8").Value = "  -1.31%  "
Set-TextValue "D9" "0.505"
$ws.Range("E9").Value = "  +0.89%  "
Set-TextValue "D10" "0.140"
$ws.Range("E10").Value = "  +0.16%  "
Set-TextValue "D11" "5.31"
$ws.Range("E11").Value = "  +1.78%  "
Set-TextValue "D12" "0.447"
$ws.Range("E12").Value = "  -0.46%  "
Set-TextValue "D13" "0.0000238"
$ws.Range("E13").Value = "  +2.17%  "
Set-TextValue "D14" "34.13"
$ws.Range("E14").Value = "  +3.75%  "
Set-TextValue "D15" "3.658.26"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("E16").Value = "  +1.34%  "
Set-TextValue "D17" "3.144.10"
$ws.Range("E17").Value = "  -0.98%  "
Set-TextValue "D18" "63.007.45"
$ws.Range("E18").Value = "  +0.96%  "
Set-TextValue "D19" "6.46"
$ws.Range("E19").Value = "  -0.89%  "
Set-TextValue "D20" "454.65"
$ws.Range("E20").Value = "  -0.03%  "
Set-TextValue "D21" "13.74"
$ws.Range("E21").Value = "  -0.13%  "
Set-TextValue "D22" "0.685"
$ws.Range("E22").Value = "  -1.92%  "
Set-TextValue "D23" "7.50"
$ws.Range("E23").Value = "  -1.09%  "
Set-TextValue "D24" "13.05"
$ws.Range("E24").Value = "  -1.92%  "
Set-TextValue "D25" "82.02"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +0.08%  "
Set-TextValue "D28" "2.64"
$ws.Range("E28").Value = "  -0.84%  "
Set-TextValue "D29" "2.05"
$ws.Range("E29").Value = "  +2.41%  "
Set-TextValue "D30" "7.56"
$ws.Range("E30").Value = "  -2.57%  "
Set-TextValue "D31" "6.59"
$ws.Range("E31").Value = "  -4.44%  "
Set-TextValue "D32" "26.62"
$ws.Range("E32").Value = "  -1.87%  "
Set-TextValue "D33" "0.0993"
$ws.Range("E33").Value = "  -1.92%  "
Set-TextValue "D34" "2.37"
$ws.Range("E34").Value = "  -0.47%  "
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  -2.44%  "
Set-TextValue "D36" "5.80"
$ws.Range("E36").Value = "  +0.62%  "
Set-TextValue "D37" "50.70"
$ws.Range("E37").Value = "  -0.71%  "
Set-TextValue "D38" "0.0₃0713"
$ws.Range("E38").Value = "  +3.82%  "
Set-TextValue "D39" "0.0383"
$ws.Range("E39").Value = "  +0.50%  "
Set-TextValue "D40" "8.03"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  -1.12%  "
Set-TextValue "D42" "2.57"
$ws.Range("E42").Value = "  -1.59%  "
Set-TextValue "D43" "383.10"
$ws.Range("E43").Value = "  -6.37%  "
Set-TextValue "D44" "2.755.53"
$ws.Range("E44").Value = "  -5.83%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D45" "0.999"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D46" "0.246"
$ws.Range("E46").Value = "  -0.92%  "
Set-TextValue "D47" "125.72"
$ws.Range("E47").Value = "  +0.80%  "
Set-TextValue "D48" "35.00"
$ws.Range("E48").Value = "  -1.23%  "
Set-TextValue "D49" "2.07"
$ws.Range("E49").Value = "  -2.36%  "
Set-TextValue "D50" "0.110"
$ws.Range("E50").Value = "  +0.29%  "
Set-TextValue "D51" "24.54"
$ws.Range("E51").Value = "  -2.46%  "
